$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, whether the value must be forced
# to Text so Excel does not auto-convert a numeric-looking string (prices
# like "236.06" or "1.960") into a real number and lose trailing zeros /
# thousands-dot formatting.
$updates = @(
    @{ Cell = 'D2'; Value = '30.295.24'; ForceText = $true }
    @{ Cell = 'E2'; Value = '  -0.17%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '1.868.53'; ForceText = $true }
    @{ Cell = 'E3'; Value = '  +0.04%  '; ForceText = $false }
    @{ Cell = 'E4'; Value = '  +0.01%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '236.06'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  +0.02%  '; ForceText = $false }
    @{ Cell = 'D7'; Value = '0.4703'; ForceText = $true }
    @{ Cell = 'E7'; Value = '  +0.54%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '0.2905'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  +2.10%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '0.06621'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  +1.24%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '21.64'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  -1.39%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '0.08013'; ForceText = $true }
    @{ Cell = 'D12'; Value = '97.34'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  -0.32%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '1.873.83'; ForceText = $true }
    @{ Cell = 'E13'; Value = '  -0.06%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '5.148'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  -0.23%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '0.6858'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  +0.81%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '273.72'; ForceText = $true }
    @{ Cell = 'E16'; Value = '  -2.49%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '30.299.70'; ForceText = $true }
    @{ Cell = 'E17'; Value = '  -0.15%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '14.05'; ForceText = $true }
    @{ Cell = 'E18'; Value = '  +6.18%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '0.000007709'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  +5.28%  '; ForceText = $false }
    @{ Cell = 'E20'; Value = '  +0.05%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '2.117.48'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  +0.04%  '; ForceText = $false }
    @{ Cell = 'B22'; Value = 'Uniswap'; ForceText = $false }
    @{ Cell = 'C22'; Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; ForceText = $false }
    @{ Cell = 'D22'; Value = '5.301'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  -2.19%  '; ForceText = $false }
    @{ Cell = 'B23'; Value = 'BinanceUSD'; ForceText = $false }
    @{ Cell = 'C23'; Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; ForceText = $false }
    @{ Cell = 'D23'; Value = '1.001'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  -0.03%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '6.208'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  +0.75%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '167.81'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  +0.95%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '9.266'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  +0.97%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '18.97'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  -0.76%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '1.958'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  +0.87%  '; ForceText = $false }
    @{ Cell = 'D29'; Value = '1.368'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  -1.61%  '; ForceText = $false }
    @{ Cell = 'D30'; Value = '0.09951'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  +1.73%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '4.362'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  -1.08%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '1.463'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  -1.24%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '4.083'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  -0.54%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '0.04708'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  -0.78%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '1.131'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  -0.49%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '0.7039'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  -0.95%  '; ForceText = $false }
    @{ Cell = 'D37'; Value = '2.707'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  -0.35%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '0.01879'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  +0.56%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '2.637'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  +3.00%  '; ForceText = $false }
    @{ Cell = 'D40'; Value = '6.331'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  -0.27%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '73.31'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  -2.08%  '; ForceText = $false }
    @{ Cell = 'D42'; Value = '1.960'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  -0.88%  '; ForceText = $false }
    @{ Cell = 'D43'; Value = '0.8407'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  -1.52%  '; ForceText = $false }
    @{ Cell = 'D44'; Value = '0.4163'; ForceText = $true }
    @{ Cell = 'E44'; Value = '  -0.88%  '; ForceText = $false }
    @{ Cell = 'B45'; Value = 'PaxDollar'; ForceText = $false }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; ForceText = $false }
    @{ Cell = 'D45'; Value = '0.9998'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  -0.02%  '; ForceText = $false }
    @{ Cell = 'B46'; Value = 'Quant'; ForceText = $false }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; ForceText = $false }
    @{ Cell = 'D46'; Value = '103.76'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  +0.21%  '; ForceText = $false }
    @{ Cell = 'B47'; Value = 'EnergySwap'; ForceText = $false }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; ForceText = $false }
    @{ Cell = 'D47'; Value = '9.297'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  -1.33%  '; ForceText = $false }
    @{ Cell = 'B48'; Value = 'Aptos'; ForceText = $false }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; ForceText = $false }
    @{ Cell = 'D48'; Value = '7.097'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  -1.88%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '936.38'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  -3.72%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '34.44'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  +0.67%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '0.05668'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  +0.43%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
